# The SNCH pharmacy-report test pack no longer includes a "Supplier Stock
# Report" test case, so remove its row (TC010 /
# Pharmacy\Reports\Stock\TC01SupplierStockReport.py) from the sheet. Deleting
# the entire row shifts the subsequent rows up and lets Excel drop the now
# unused shared strings automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(9).Delete()

# Match the author's final cursor position in the sheet.
$ws.Range("C14").Select() | Out-Null
